$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (rows 1-10): fill in the new trailing values, a new "2,NNNNN"
# series (step 5).
$dValues = @{
    1  = "2,37851"
    2  = "2,37856"
    3  = "2,37861"
    4  = "2,37866"
    5  = "2,37871"
    6  = "2,37876"
    7  = "2,37881"
    8  = "2,37886"
    9  = "2,37891"
    10 = "2,37896"
}
foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# Column C (rows 11-31): fill in the new trailing values, continuing the
# existing "2,NNNNN" series (step 4).
$cValues = @{
    11 = "2,37766"
    12 = "2,37770"
    13 = "2,37774"
    14 = "2,37778"
    15 = "2,37782"
    16 = "2,37786"
    17 = "2,37790"
    18 = "2,37794"
    19 = "2,37798"
    20 = "2,37802"
    21 = "2,37806"
    22 = "2,37810"
    23 = "2,37814"
    24 = "2,37818"
    25 = "2,37822"
    26 = "2,37826"
    27 = "2,37830"
    28 = "2,37834"
    29 = "2,37838"
    30 = "2,37842"
    31 = "2,37846"
}
foreach ($row in $cValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $cValues[$row]
}

# Move the selection from L32 to G8.
$ws.Range("G8").Select()
